# Applies a weekly data update to the "Jengibre" sheet:
# - A brand new observation is inserted as the first record (row 250), and all
#   the existing records from row 250 through 362 shift down by one row,
#   with the last existing record (old row 362) becoming new row 363.
# Only the columns that vary per-record (D=Fecha, I=Calidad, J=Volumen,
# K=Precio minimo, L=Precio maximo, M=Precio promedio ponderado, P=Precio $/Kg)
# need to move; the remaining columns (A,B,C,E,F,G,H,N,O,Q,R) are identical
# for every record in this sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 250
$lastDataRow = 362
$newLastRow = 363

# --- Step 1: capture the current per-record values for rows 250..362 ---
$D = @{}
$I = @{}
$J = @{}
$K = @{}
$L = @{}
$M = @{}
$P = @{}

for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $D[$r] = $ws.Cells.Item($r, 4).Value2()
    $I[$r] = $ws.Cells.Item($r, 9).Value2()
    $J[$r] = $ws.Cells.Item($r, 10).Value2()
    $K[$r] = $ws.Cells.Item($r, 11).Value2()
    $L[$r] = $ws.Cells.Item($r, 12).Value2()
    $M[$r] = $ws.Cells.Item($r, 13).Value2()
    $P[$r] = $ws.Cells.Item($r, 16).Value2()
}

# --- Step 2: create the new last row (363) as a full copy of old row 362 ---
for ($c = 1; $c -le 18; $c++) {
    $ws.Cells.Item($newLastRow, $c).Value = $ws.Cells.Item($lastDataRow, $c).Value2()
}
# preserve the date number format on column D (Fecha)
$ws.Cells.Item($newLastRow, 4).NumberFormat = $ws.Cells.Item($lastDataRow, 4).NumberFormat()

# --- Step 3: shift rows 362 down to 251, each one takes the previous row's values ---
for ($r = $lastDataRow; $r -ge ($firstDataRow + 1); $r--) {
    $src = $r - 1
    $ws.Cells.Item($r, 4).Value = $D[$src]
    $ws.Cells.Item($r, 9).Value = $I[$src]
    $ws.Cells.Item($r, 10).Value = $J[$src]
    $ws.Cells.Item($r, 11).Value = $K[$src]
    $ws.Cells.Item($r, 12).Value = $L[$src]
    $ws.Cells.Item($r, 13).Value = $M[$src]
    $ws.Cells.Item($r, 16).Value = $P[$src]
}

# --- Step 4: set the brand new first record values on row 250 ---
$ws.Cells.Item($firstDataRow, 4).Value = 45202
$ws.Cells.Item($firstDataRow, 9).Value = "Primera"
$ws.Cells.Item($firstDataRow, 10).Value = 35
$ws.Cells.Item($firstDataRow, 11).Value = 24000
$ws.Cells.Item($firstDataRow, 12).Value = 24000
$ws.Cells.Item($firstDataRow, 13).Value = 24000
$ws.Cells.Item($firstDataRow, 16).Value = 1846
